# Minor fix for test.xlsx.
#
# Summary of the change:
#   - Sheet1 gains a new header cell D1 = "header4" (styled like the other
#     header cells A1:C1).
#   - Sheet1 becomes the active/selected sheet (it was "2ndsheet" before),
#     with the active cell/selection set to D2.
#   - "2ndsheet" is no longer the selected sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet1: add the new "header4" column header ------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("D1").Value = "header4"
# Match the bold header formatting already used by A1:C1.
$ws1.Range("D1").Font.Bold = $true

# --- Make Sheet1 the active sheet, with D2 selected ----------------------
$ws1.Activate() | Out-Null
$ws1.Range("D2").Select() | Out-Null
